$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058167944988267
$ws.Range("D2").Value = 1.058725373056155
$ws.Range("E2").Value = 1.062588550130235
$ws.Range("F2").Value = 1.071747681098949
$ws.Range("I2").Value = 1.051361620810624
$ws.Range("J2").Value = 1.063160570128734
$ws.Range("K2").Value = 1.061456822395036
$ws.Range("L2").Value = 1.065309489731029
$ws.Range("M2").Value = 1.074444030770481
$ws.Range("N2").Value = 1.064670380084539

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059510692235543
$ws.Range("D3").Value = 1.059793389953043
$ws.Range("E3").Value = 1.064148733673818
$ws.Range("F3").Value = 1.073410820050889
$ws.Range("I3").Value = 1.051846399970022
$ws.Range("J3").Value = 1.064154186863574
$ws.Range("K3").Value = 1.06233828545679
$ws.Range("L3").Value = 1.06668264594684
$ws.Range("M3").Value = 1.075921694277832
$ws.Range("N3").Value = 1.065665407869112

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060377981668332
$ws.Range("D4").Value = 1.060483064905856
$ws.Range("E4").Value = 1.065154928149513
$ws.Range("F4").Value = 1.074484279259166
$ws.Range("I4").Value = 1.052157958495326
$ws.Range("J4").Value = 1.064795083169943
$ws.Range("K4").Value = 1.06290665139731
$ws.Range("L4").Value = 1.067567332066591
$ws.Range("M4").Value = 1.076874664656871
$ws.Range("N4").Value = 1.066307214321747

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060742224202449
$ws.Range("D5").Value = 1.060772673207596
$ws.Range("E5").Value = 1.065577144230937
$ws.Range("F5").Value = 1.074934927791413
$ws.Range("I5").Value = 1.052288431762263
$ws.Range("J5").Value = 1.065064033105569
$ws.Range("K5").Value = 1.063145118487129
$ws.Range("L5").Value = 1.06793834798816
$ws.Range("M5").Value = 1.077274545282428
$ws.Range("N5").Value = 1.066576546197131

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060803360870062
$ws.Range("D6").Value = 1.060821280445305
$ws.Range("E6").Value = 1.065647990348522
$ws.Range("F6").Value = 1.07501055688773
$ws.Range("I6").Value = 1.052310309217132
$ws.Range("J6").Value = 1.065109162827697
$ws.Range("K6").Value = 1.063185130505518
$ws.Range("L6").Value = 1.068000590424118
$ws.Range("M6").Value = 1.077341643487088
$ws.Range("N6").Value = 1.066621740008641

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060382850124171
$ws.Range("D7").Value = 1.060486935960893
$ws.Range("E7").Value = 1.065160572901314
$ws.Range("F7").Value = 1.074490303320516
$ws.Range("I7").Value = 1.052159703869769
$ws.Range("J7").Value = 1.064798678782679
$ws.Range("K7").Value = 1.062909839662696
$ws.Range("L7").Value = 1.067572293144365
$ws.Range("M7").Value = 1.076880010805703
$ws.Range("N7").Value = 1.066310815040666

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058622059495344
$ws.Range("D8").Value = 1.059086607928271
$ws.Range("E8").Value = 1.063116523021613
$ws.Range("F8").Value = 1.072310314420404
$ws.Range("I8").Value = 1.051525896234756
$ws.Range("J8").Value = 1.06349679333821
$ws.Range("K8").Value = 1.061755133870346
$ws.Range("L8").Value = 1.065774356579685
$ws.Range("M8").Value = 1.07494407931718
$ws.Range("N8").Value = 1.065007080769538

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.05550707485115
$ws.Range("D9").Value = 1.056608080192013
$ws.Range("E9").Value = 1.059488388096098
$ws.Range("F9").Value = 1.068447593005431
$ws.Range("I9").Value = 1.050392626597349
$ws.Range("J9").Value = 1.061186829430176
$ws.Range("K9").Value = 1.059704865461536
$ws.Range("L9").Value = 1.062576200568616
$ws.Range("M9").Value = 1.071507833037924
$ws.Range("N9").Value = 1.062693836447798

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053421740212671
$ws.Range("D10").Value = 1.05494803782139
$ws.Range("E10").Value = 1.057051119146853
$ws.Range("F10").Value = 1.065857278994093
$ws.Range("I10").Value = 1.049625880509076
$ws.Range("J10").Value = 1.059635826981102
$ws.Range("K10").Value = 1.058327275796743
$ws.Range("L10").Value = 1.060423152993143
$ws.Range("M10").Value = 1.069199486423959
$ws.Range("N10").Value = 1.061140631397343

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052516609152738
$ws.Range("D11").Value = 1.0542273261195
$ws.Range("E11").Value = 1.055991172175288
$ws.Range("F11").Value = 1.06473185858103
$ws.Range("I11").Value = 1.049291163056004
$ws.Range("J11").Value = 1.058961531088327
$ws.Range("K11").Value = 1.057728147630135
$ws.Range("L11").Value = 1.059485713583325
$ws.Range("M11").Value = 1.068195617394537
$ws.Range("N11").Value = 1.060465377927058

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.052180068851876
$ws.Range("D12").Value = 1.053959329674816
$ws.Range("E12").Value = 1.055596753961426
$ws.Range("F12").Value = 1.064313240493203
$ws.Range("I12").Value = 1.049166422909274
$ws.Range("J12").Value = 1.058710654749415
$ws.Range("K12").Value = 1.05750520473352
$ws.Range("L12").Value = 1.059136716808919
$ws.Range("M12").Value = 1.067822068588122
$ws.Range("N12").Value = 1.060214145314968

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052252273158259
$ws.Range("D13").Value = 1.054016829137383
$ws.Range("E13").Value = 1.055681390255889
$ws.Range("F13").Value = 1.064403062313613
$ws.Range("I13").Value = 1.049193198762998
$ws.Range("J13").Value = 1.058764487368922
$ws.Range("K13").Value = 1.057553044940075
$ws.Range("L13").Value = 1.059211613802627
$ws.Range("M13").Value = 1.06790222651018
$ws.Range("N13").Value = 1.06026805438297

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.052488797497214
$ws.Range("D14").Value = 1.05420517944328
$ws.Range("E14").Value = 1.055958583956106
$ws.Range("F14").Value = 1.064697267545647
$ws.Range("I14").Value = 1.049280860406861
$ws.Range("J14").Value = 1.05894080203371
$ws.Range("K14").Value = 1.057709727295856
$ws.Range("L14").Value = 1.059456881607659
$ws.Range("M14").Value = 1.068164753431602
$ws.Range("N14").Value = 1.060444619434806

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052634483385871
$ws.Range("D15").Value = 1.054321189399314
$ws.Range("E15").Value = 1.05612927799721
$ws.Range("F15").Value = 1.064878458846179
$ws.Range("I15").Value = 1.049334817037933
$ws.Range("J15").Value = 1.059049380406195
$ws.Range("K15").Value = 1.057806211279946
$ws.Range("L15").Value = 1.059607894044059
$ws.Range("M15").Value = 1.068326416072171
$ws.Range("N15").Value = 1.060553352001034

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.053481764375163
$ws.Range("D16").Value = 1.054995828473169
$ws.Range("E16").Value = 1.05712136614487
$ws.Range("F16").Value = 1.065931888044703
$ws.Range("I16").Value = 1.049648037153032
$ws.Range("J16").Value = 1.059680520268793
$ws.Range("K16").Value = 1.058366982192631
$ws.Range("L16").Value = 1.060485257820171
$ws.Range("M16").Value = 1.069266017182148
$ws.Range("N16").Value = 1.061185388154629

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.054012654967341
$ws.Range("D17").Value = 1.055418498047836
$ws.Range("E17").Value = 1.057742435368789
$ws.Range("F17").Value = 1.066591648215651
$ws.Range("I17").Value = 1.049843783260423
$ws.Range("J17").Value = 1.06007568946132
$ws.Range("K17").Value = 1.058718032583795
$ws.Range("L17").Value = 1.06103421376621
$ws.Range("M17").Value = 1.069854231616669
$ws.Range("N17").Value = 1.061581118532737

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.054322105937566
$ws.Range("D18").Value = 1.055664850929679
$ws.Range("E18").Value = 1.058104252469718
$ws.Range("F18").Value = 1.066976109882881
$ws.Range("I18").Value = 1.049957697309803
$ws.Range("J18").Value = 1.060305924773474
$ws.Range("K18").Value = 1.058922541688058
$ws.Range("L18").Value = 1.061353914093728
$ws.Range("M18").Value = 1.070196909456857
$ws.Range("N18").Value = 1.061811680805442

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054427585612692
$ws.Range("D19").Value = 1.055748820020722
$ws.Range("E19").Value = 1.058227548268029
$ws.Range("F19").Value = 1.067107139909977
$ws.Range("I19").Value = 1.049996494843696
$ws.Range("J19").Value = 1.06038438514055
$ws.Range("K19").Value = 1.058992231333607
$ws.Range("L19").Value = 1.061462839913165
$ws.Range("M19").Value = 1.070313683353983
$ws.Range("N19").Value = 1.061890252595239

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.053955717065033
$ws.Range("D20").Value = 1.0553731685909
$ws.Range("E20").Value = 1.057675846366836
$ws.Range("F20").Value = 1.066520900096363
$ws.Range("I20").Value = 1.049822808616143
$ws.Range("J20").Value = 1.060033318511296
$ws.Range("K20").Value = 1.058680394376391
$ws.Range("L20").Value = 1.060975367446539
$ws.Range("M20").Value = 1.06979116509013
$ws.Range("N20").Value = 1.061538687411103

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.052419156273868
$ws.Range("D21").Value = 1.054149723084385
$ws.Range("E21").Value = 1.055876976963229
$ws.Range("F21").Value = 1.064610647769218
$ws.Range("I21").Value = 1.049255057619263
$ws.Range("J21").Value = 1.058888893211422
$ws.Range("K21").Value = 1.05766359933813
$ws.Range("L21").Value = 1.059384678286136
$ws.Range("M21").Value = 1.068087464335713
$ws.Range("N21").Value = 1.060392636896036

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051451120851826
$ws.Range("D22").Value = 1.053378801618636
$ws.Range("E22").Value = 1.054741859194127
$ws.Range("F22").Value = 1.063406192848196
$ws.Range("I22").Value = 1.048895709136673
$ws.Range("J22").Value = 1.058166955657894
$ws.Range("K22").Value = 1.057021982543531
$ws.Range("L22").Value = 1.058379970185533
$ws.Range("M22").Value = 1.067012412499037
$ws.Range("N22").Value = 1.059669674108375

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.051964481433692
$ws.Range("D23").Value = 1.053787644216524
$ws.Range("E23").Value = 1.055344000553052
$ws.Range("F23").Value = 1.06404502529432
$ws.Range("I23").Value = 1.049086433538256
$ws.Range("J23").Value = 1.058549897477992
$ws.Range("K23").Value = 1.057362337421918
$ws.Range("L23").Value = 1.058913024171119
$ws.Range("M23").Value = 1.067582689518464
$ws.Range("N23").Value = 1.060053159749782

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053981445495796
$ws.Range("D24").Value = 1.055393651586638
$ws.Range("E24").Value = 1.057705936434745
$ws.Range("F24").Value = 1.066552869253148
$ws.Range("I24").Value = 1.049832286962219
$ws.Range("J24").Value = 1.060052464923113
$ws.Range("K24").Value = 1.058697402242325
$ws.Range("L24").Value = 1.061001959097542
$ws.Range("M24").Value = 1.069819663427846
$ws.Range("N24").Value = 1.061557861013022

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.056313869812652
$ws.Range("D25").Value = 1.057250171274857
$ws.Range("E25").Value = 1.060429547466863
$ws.Range("F25").Value = 1.069448809915755
$ws.Range("I25").Value = 1.050687569553421
$ws.Range("J25").Value = 1.061785927884796
$ws.Range("K25").Value = 1.06023678028977
$ws.Range("L25").Value = 1.063406632985007
$ws.Range("M25").Value = 1.072399214509136
$ws.Range("N25").Value = 1.06329378569094

